# Fixing DATA_TYPE / DATA_LENGTH for both Oracle (DEPT) and MySQL (EMP) metadata:
# every "decimal" / "decimal(x,y)" entry becomes "int".

$wb = $excel.ActiveWorkbook

$empSheet  = $wb.Worksheets.Item("EMP")
$deptSheet = $wb.Worksheets.Item("DEPT")

# EMP sheet: empno, mgr, sal, comm, deptno were decimal/decimal(n,m) -> int
$empSheet.Range("B2").Value = "int"
$empSheet.Range("C2").Value = "int"

$empSheet.Range("B5").Value = "int"
$empSheet.Range("C5").Value = "int"

$empSheet.Range("B7").Value = "int"
$empSheet.Range("C7").Value = "int"

$empSheet.Range("B8").Value = "int"
$empSheet.Range("C8").Value = "int"

$empSheet.Range("B9").Value = "int"
$empSheet.Range("C9").Value = "int"

# DEPT sheet: deptno was decimal/decimal(2,0) -> int
$deptSheet.Range("B2").Value = "int"
$deptSheet.Range("C2").Value = "int"

# Update the last-used selection on each sheet to match the author's saved state.
# Select DEPT first, then EMP last, so EMP remains the active/displayed tab
# (matches tabSelected="1" staying on the EMP sheet in the saved file).
$deptSheet.Range("C14").Select()
$empSheet.Range("E9").Select()
